# Update "想去人数" (want-to-go count) figures across the workbook's sheets
# to match freshly re-scraped output, per commit "Update gh-pages to output
# generated at 7921097".

$wb = $excel.ActiveWorkbook

# Sheet "展览" (Exhibitions)
$ws = $wb.Worksheets.Item("展览")
$ws.Range("F5").Value  = 146
$ws.Range("F6").Value  = 570
$ws.Range("F7").Value  = 261
$ws.Range("F9").Value  = 992
$ws.Range("F11").Value = 2220
$ws.Range("F12").Value = 707
$ws.Range("F15").Value = 786
$ws.Range("F20").Value = 11
$ws.Range("F24").Value = 1242
$ws.Range("F30").Value = 508
$ws.Range("F33").Value = 259
$ws.Range("F38").Value = 888
$ws.Range("F40").Value = 62
$ws.Range("F41").Value = 149
$ws.Range("F43").Value = 204

# Sheet "演出" (Performances)
$ws = $wb.Worksheets.Item("演出")
$ws.Range("F8").Value  = 659
$ws.Range("F19").Value = 614
$ws.Range("F22").Value = 437
$ws.Range("F25").Value = 186

# Sheet "本地生活" (Local life)
$ws = $wb.Worksheets.Item("本地生活")
$ws.Range("F6").Value  = 2154
$ws.Range("F11").Value = 853
$ws.Range("F12").Value = 140

# Sheet "全部类型" (All types, aggregated list)
$ws = $wb.Worksheets.Item("全部类型")
$ws.Range("F5").Value  = 2154
$ws.Range("F12").Value = 146
$ws.Range("F14").Value = 570
$ws.Range("F15").Value = 261
$ws.Range("F17").Value = 992
$ws.Range("F18").Value = 853
$ws.Range("F21").Value = 140
$ws.Range("F23").Value = 659
$ws.Range("F24").Value = 786
$ws.Range("F31").Value = 1242
$ws.Range("F36").Value = 508
$ws.Range("F38").Value = 259
$ws.Range("F42").Value = 888
$ws.Range("F46").Value = 186
$ws.Range("F48").Value = 149
$ws.Range("F50").Value = 204
